$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B5").Value = 7.3973141682915102
$ws.Range("C5").Value = 0.14814375886436101
$ws.Range("D5").Value = 7.6946542968670402

$ws.Range("D7").Value = 22.813350029660832

$ws.Range("B8").Value = 27.465494818840153
$ws.Range("C8").Value = 1.6260938337383766
$ws.Range("D8").Value = 9.9562596473047087

$ws.Range("B11").Value = 36.816468749811449
$ws.Range("D11").Value = 20.882440305774484

$ws.Range("B12").Value = 41.523553938490778
$ws.Range("D12").Value = 26.885059555058994

$ws.Range("B14").Value = 19.486153847891028
$ws.Range("C14").Value = 0.29147226075451121
$ws.Range("D14").Value = 12.16107844354066

$ws.Range("B17").Value = 45.383661056296518
$ws.Range("C17").Value = 1.4182012161285864
$ws.Range("D17").Value = 29.33246448540104

$ws.Range("B18").Value = 47.718969729854848
$ws.Range("C18").Value = 0.46524782763105599
$ws.Range("D18").Value = 7.1924510678553126

$ws.Range("B19").Value = 13.492639850759762
$ws.Range("D19").Value = 43.426468720724202

$ws.Range("B20").Value = 38.934076480058238
$ws.Range("C20").Value = 2.3731203348415986
$ws.Range("D20").Value = 7.3716287036456434

$ws.Range("B21").Value = 40.205527461993398
$ws.Range("D21").Value = 22.705359783795299

$ws.Range("B22").Value = 34.375830990633247
$ws.Range("D22").Value = 25.982993370849659

$ws.Range("B23").Value = 34.326143140675796
$ws.Range("D23").Value = 36.268372706706558

$ws.Range("B24").Value = 45.732591613085717
$ws.Range("D24").Value = 17.180232459795224

$ws.Range("B25").Value = 48.223490610689609
$ws.Range("D25").Value = 13.106567106487848

$ws.Range("B26").Value = 27.881907098690359
$ws.Range("C26").Value = 1.4530734990347851
$ws.Range("D26").Value = 3.9241209721414494

$ws.Range("B27").Value = 36.772546695502257
$ws.Range("C27").Value = 1.1535423607123096
$ws.Range("D27").Value = 2.6869815292431616

$ws.Range("B28").Value = 34.472520598432808
$ws.Range("C28").Value = 1.0297316322904859
$ws.Range("D28").Value = 10.01531539587117
